# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Septiembre de 2020 a las 21:22"

# --- Swap display order of "Santa Lucia" / "Timor Oriental" ---
# Row 207 previously showed "Santa Lucia", row 208 previously showed "Timor Oriental".
# Their statistics are identical, only the shared-string ordering changes, so the
# country names in those two rows must be swapped.
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("A208").Value = "Santa Lucia"

# --- Update statistics for Estados Unidos (row 4) ---
$ws.Range("B4").Value = 7309027
$ws.Range("C4").Value = 21466
$ws.Range("D4").Value = 4537020
$ws.Range("E4").Value = 2562676
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 154
$ws.Range("H4").Value = 209331

# --- Update statistics for India (row 5) ---
$ws.Range("B5").Value = 6053010
$ws.Range("C5").Value = 62429
$ws.Range("D5").Value = 4998519
$ws.Range("E5").Value = 959329
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 628
$ws.Range("H5").Value = 95162

# --- Update statistics for Alemania (row 25) ---
$ws.Range("B25").Value = 286159
$ws.Range("C25").Value = 1134
$ws.Range("D25").Value = 250800
$ws.Range("E25").Value = 25825

# --- Update statistics for Ecuador (row 30) ---
$ws.Range("B30").Value = 134747
$ws.Range("C30").Value = 766
$ws.Range("D30").Value = 112296
$ws.Range("E30").Value = 11172
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = 11279
